# daily auto push: 2025-10-07 18:39 UTC
# Append one new log row (row 77) to the sheet: 2025/10/08, 水, 2, 201
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 77

# Column A holds a date-like string ("2025/10/08") that must stay plain
# text, not get auto-converted into a date serial number. Mark the cell
# as Text before writing it, then drop the formatting again so the new
# row doesn't pick up a stray number format/style (matching the rest of
# the sheet's unstyled data rows).
$ws.Range("A" + $row).NumberFormat = "@"
$ws.Range("A" + $row).Value = "2025/10/08"
$ws.Range("A" + $row).ClearFormats()

$ws.Range("B" + $row).Value = "水"
$ws.Range("C" + $row).Value = 2
$ws.Range("D" + $row).Value = 201
